$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new officer netid for the Cottage club (column F, row 2)
$ws.Range("F2").Value = "ad15"
